$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.350.12'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.634.32'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.76'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.19'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.632.93'
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.359'
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.115.99'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '72.251.07'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.73'
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.631.99'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.05'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '376.05'
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.86'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.50'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.24'
$ws.Range("E26").Value = '  -2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.41'
$ws.Range("E27").Value = '  -3.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.770.27'
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0948'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '490.86'
$ws.Range("E32").Value = '  -3.78%  '
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.50'
$ws.Range("E37").Value = '  +8.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.15'
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  -4.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  -2.55%  '
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.21'
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.542'
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("E51").Value = '  +0.81%  '
